{"js": "// The document currently has several paragraphs of placeholder/test text.\n// The edit collapses all of that into a single paragraph holding new text,\n// while keeping the final paragraph's formatting and its \"_GoBack\" bookmark.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Delete every paragraph except the last one; deleting a Word.Paragraph\n// removes the paragraph (and its mark) entirely, leaving the following\n// paragraph (here, the last one) with its own formatting/bookmarks intact.\nconst count = paragraphs.items.length;\nfor (let i = 0; i < count - 1; i++) {\n  paragraphs.items[i].delete();\n}\nawait context.sync();\n\n// Re-fetch the (now single) remaining paragraph and replace its text.\nconst remaining = context.document.body.paragraphs;\nremaining.load(\"items\");\nawait context.sync();\n\nconst target = remaining.items[0];\ntarget.insertText(\n  \"20171128 \\u2013 Let\\u2019s upload this to jerrycon/documents on github to see what happens to the previously existing, identically named document.\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Merge every paragraph in the body into the last paragraph by repeatedly\n# deleting the paragraph mark that ends the first remaining paragraph.\n# This keeps the last paragraph's formatting/bookmarks (matching real Word\n# \"delete a paragraph mark merges forward into the following paragraph\").\nwhile ($d.Paragraphs.Count -gt 1) {\n    $p = $d.Paragraphs.Item(1)\n    $markStart = $p.Range.End - 1\n    $markEnd = $p.Range.End\n    $mark = $d.Range($markStart, $markEnd)\n    $mark.Delete()\n}\n\n# Replace the remaining paragraph's text (but not its paragraph mark) with\n# the new content.\n$p1 = $d.Paragraphs.Item(1)\n$textRange = $d.Range($p1.Range.Start, $p1.Range.End - 1)\n$textRange.Text = \"20171128 \u2013 Let\u2019s upload this to jerrycon/documents on github to see what happens to the previously existing, identically named document.\"\n"}
